$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2017-02-22 07:16:56"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2017-02-22 07:16:56"
$wsDeDe.Range("L4").Value = "2017-02-22 07:18:01"
